$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Target cluster" value in column D for rows 2-6 from "Resolving-Mac" to "Inflammatory-Mac"
$ws.Range("D2:D6").Value = "Inflammatory-Mac"

# Row 2
$ws.Range("G2").Value = 0.9024946666666667
$ws.Range("H2").Value = 2.707484
$ws.Range("I2").Value = 0.01048932664802141
$ws.Range("J2").Value = 0.01058096843676313
$ws.Range("M2").Value = 0.08893566666666668
$ws.Range("N2").Value = 0.266807
$ws.Range("Q2").Value = 0.08026396484311112
$ws.Range("R2").Value = 0.7223756835880001
$ws.Range("S2").Value = 0.01048932664802141
$ws.Range("T2").Value = 0.01058096843676313

# Row 3
$ws.Range("H3").Value = 4.840617
$ws.Range("I3").Value = 0.01875350432023437
$ws.Range("J3").Value = 0.01891734750471621
$ws.Range("M3").Value = 0.08893566666666668
$ws.Range("N3").Value = 0.266807
$ws.Range("Q3").Value = 0.1435011666576667
$ws.Range("R3").Value = 1.291510499919
$ws.Range("S3").Value = 0.01875350432023437
$ws.Range("T3").Value = 0.01891734750471621

# Row 4
$ws.Range("G4").Value = 36.42588166666667
$ws.Range("H4").Value = 109.277645
$ws.Range("I4").Value = 0.4233631348261055
$ws.Range("J4").Value = 0.4270619189582678
$ws.Range("M4").Value = 0.08893566666666668
$ws.Range("N4").Value = 0.266807
$ws.Range("Q4").Value = 3.239560069946112
$ws.Range("R4").Value = 29.156040629515
$ws.Range("S4").Value = 0.4233631348261055
$ws.Range("T4").Value = 0.4270619189582678

# Row 5
$ws.Range("G5").Value = 2.2355605
$ws.Range("H5").Value = 4.471121
$ws.Range("I5").Value = 0.02598300598553572
$ws.Range("J5").Value = 0.01747334062840217
$ws.Range("M5").Value = 0.08893566666666668
$ws.Range("N5").Value = 0.266807
$ws.Range("Q5").Value = 0.1988210634411667
$ws.Range("R5").Value = 1.192926380647
$ws.Range("S5").Value = 0.02598300598553572
$ws.Range("T5").Value = 0.01747334062840217

# Row 6
$ws.Range("G6").Value = 44.86185700000001
$ws.Range("H6").Value = 134.585571
$ws.Range("I6").Value = 0.5214110282201031
$ws.Range("J6").Value = 0.5259664244718506
$ws.Range("M6").Value = 0.08893566666666668
$ws.Range("N6").Value = 0.266807
$ws.Range("Q6").Value = 3.989819160199668
$ws.Range("R6").Value = 35.90837244179701
$ws.Range("S6").Value = 0.5214110282201031
$ws.Range("T6").Value = 0.5259664244718506
